$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price record (week of 2021-10-15, serial 44491) needs to be
# inserted for "Provincia del Elquí" ahead of the existing row 116 record
# (serial 44301). Insert a fresh row at 116, pushing the existing rows
# 116-133 down to 117-134, then populate the new row with its own data.
$ws.Rows.Item(116).Insert()

$ws.Range("A116").Value = 7
$ws.Range("B116").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C116").Value = "Ñuble"
$ws.Range("D116").Value = 44491
$ws.Range("E116").Value = 16
$ws.Range("F116").Value = 100112017
$ws.Range("G116").Value = "Apio"
$ws.Range("H116").Value = "Americana (o)"
$ws.Range("I116").Value = "Primera"
$ws.Range("J116").Value = 120
$ws.Range("K116").Value = 8000
$ws.Range("L116").Value = 8500
$ws.Range("M116").Value = 8250
$ws.Range("N116").Value = "`$/docena de matas"
$ws.Range("O116").Value = "Provincia del Elquí"
$ws.Range("P116").Value = 1375
$ws.Range("Q116").Value = 6
$ws.Range("R116").Value = "Hortaliza"
